$wb = $excel.ActiveWorkbook

# Rename sheet3 "contingency" -> "heatmap"
$ws = $wb.Worksheets.Item(3)
$ws.Name = "heatmap"
$ws.Activate()

# Header row (row 2)
$ws.Range("B2").Value = "Excel"
$ws.Range("C2").Value = "PowerBI"
$ws.Range("D2").Value = "R"

# Row labels (column A) and data
$rows = @(
    @{ row = 3;  label = "Set up project";                                     b = 2; c = 1; d = 0 },
    @{ row = 4;  label = "See raw data for diagnosis";                         b = 2; c = 1; d = 0 },
    @{ row = 5;  label = "Factor variables";                                   b = 0; c = 1; d = 2 },
    @{ row = 6;  label = "Remove empty rows";                                  b = 0; c = 1; d = 2 },
    @{ row = 7;  label = "Filter/Select rows";                                 b = 1; c = 1; d = 2 },
    @{ row = 8;  label = "Subset columns";                                     b = 0; c = 1; d = 2 },
    @{ row = 9;  label = "Merge data";                                         b = 0; c = 1; d = 2 },
    @{ row = 10; label = "Document/Replicate data cleaning process";           b = 0; c = 1; d = 2 },
    @{ row = 11; label = "Descriptive stats, Cross-tabs, Pivot tables";        b = 2; c = 0; d = 1 },
    @{ row = 12; label = "Visualize data";                                    b = 2; c = 1; d = 1 },
    @{ row = 13; label = "Run statistical tests";                             b = 1; c = 0; d = 2 }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.row).Value = $r.label
    $ws.Range("B" + $r.row).Value = $r.b
    $ws.Range("C" + $r.row).Value = $r.c
    $ws.Range("D" + $r.row).Value = $r.d
}

# Conditional formatting: 3-color scale (red/yellow/green) on B3:D13
$cfRange = $ws.Range("B3:D13")
$cfRange.FormatConditions.AddColorScale(3) | Out-Null
$colorScale = $cfRange.FormatConditions.Item($cfRange.FormatConditions.Count)
$colorScale.ColorScaleCriteria.Item(1).Type = 1  # xlConditionValueLowestValue -> min
$colorScale.ColorScaleCriteria.Item(1).FormatColor.Color = 7039480   # RGB F8696B
$colorScale.ColorScaleCriteria.Item(2).Type = 5  # xlConditionValuePercentile
$colorScale.ColorScaleCriteria.Item(2).Value = 50
$colorScale.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167   # RGB FFEB84
$colorScale.ColorScaleCriteria.Item(3).Type = 2  # xlConditionValueHighestValue -> max
$colorScale.ColorScaleCriteria.Item(3).FormatColor.Color = 8109667   # RGB 63BE7B

# Selection on active sheet (matches the saved cursor position)
$ws.Range("H39").Select()
